$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed "CasesTab" query in B2: added missing comma after "Response to Treatment"
# and joined the trailing "Cohort" coalesce onto the RETURN clause (no blank line).
$fixedCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['MGT01'] and demo.breed in ['Australian Cattle Dog','Mixed Breed']and diag.disease_term in ['Mammary Cancer'] and diag.primary_disease_site in ['Mammary Gland']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

# New "StudyFilesTab" query in B5: file<->study association query (replaces the
# placeholder that had been duplicating the FilesTab query).
$newStudyFilesQuery = @'
 MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE 
  s.clinical_study_designation IN ['MGT01'] and
  demo.breed in ['Australian Cattle Dog','Mixed Breed'] and
  diag.disease_term in ['Mammary Cancer'] and 
  diag.primary_disease_site in ['Mammary Gland']
WITH DISTINCT f, s
RETURN 
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS `File Format`,
  coalesce(f.file_size, '') AS `Size`,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# The here-strings pick up a trailing newline from the closing "'@" line; trim it
# so the stored text matches the source query exactly (no trailing blank line).
$fixedCasesQuery = $fixedCasesQuery.TrimEnd("`r", "`n")
$newStudyFilesQuery = $newStudyFilesQuery.TrimEnd("`r", "`n")

$ws.Range("B5").Value = $newStudyFilesQuery
$ws.Range("B2").Value = $fixedCasesQuery

$ws.Rows.Item(2).RowHeight = 290
$ws.Rows.Item(5).RowHeight = 232

$ws.Range("B5").Select()
